$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$ws.Activate()

# Update the single changed cell: B5 "Meenakshi" -> "Sushil"
$ws.Range("B5").Value = "Sushil"

# Reflect the new active cell/selection on this sheet
$ws.Range("B5").Select()
